# October 5th Work.
# Fill in row 6 (2025-10-05) of the progression tracker with this day's
# entries, matching the formatting already used on the rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fills/borders) from the prior day's row (row 5) down
# onto row 6 so the new cells pick up the same "Good"/"Neutral"/"Bad" cell
# styles already used for the other completed days.
$ws.Range("B5:M5").Copy()
$ws.Range("B6:M6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Fill in the day's progress, column by column.
$ws.Range("J6").Value = "Roublard"
$ws.Range("B6").Value = "Folk Hero"
$ws.Range("C6").Value = "Esquive"
$ws.Range("E6").Value = "Baraqué"
$ws.Range("F6").Value = "Classique"
$ws.Range("G6").Value = "Athlétisme"
$ws.Range("I6").Value = "Attaque sournoise"
$ws.Range("H6").Value = "Distraction"
